$d = $word.ActiveDocument

# The document currently has a single, empty paragraph that only holds the
# "_GoBack" bookmark. We need to:
#   1. Mark the (empty) paragraph mark as Polish (pl-PL) - this produces the
#      <w:pPr><w:rPr><w:lang w:val="pl-PL"/></w:rPr></w:pPr> on the paragraph.
#   2. Insert the text "test" before the bookmark, and mark that new run as
#      Polish (pl-PL) too - this produces
#      <w:r><w:rPr><w:lang w:val="pl-PL"/></w:rPr><w:t>test</w:t></w:r>.

$para = $d.Paragraphs.Item(1)

# Step 1: set the language of the paragraph mark while the paragraph is
# still empty, so the language is stored on the paragraph's own rPr (pPr/rPr)
# rather than on a run.
$markRange = $para.Range
$markRange.LanguageID = "pl-PL"

# Step 2: insert the new text at the start of the paragraph (i.e. before the
# bookmark start/end that are already there).
$insertRange = $para.Range
$insertRange.InsertBefore("test")

# Step 3: set the language of the newly inserted "test" run.
$runRange = $d.Range(0, 4)
$runRange.LanguageID = "pl-PL"
